# Add season-record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new labels, formatted like the existing header cells.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-45: every player shares the team's overall season record.
$ws.Range("AD2:AD45").Value = 77
$ws.Range("AE2:AE45").Value = 85
$ws.Range("AF2:AF45").Value = 0
